$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.681.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.940.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.441"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.80"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.429.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.687.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.74"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.941.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "434.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.69%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.22%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.110"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.74"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.91%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.54"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.284"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.716.28"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "371.24"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.88"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.06%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.79"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.90%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.52%  "
